# Applies the "剩余" (remaining) column update across the data rows.
# For every data row (2..99) the value in column E ("剩余") is decremented
# by 1, EXCEPT:
#   - row 36, which is left untouched (its data was already stale/unchanged)
#   - row 95, which instead gets "refilled": E goes to 10 and the start
#     date in column F is reset to 2026-03-01 (20260301)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99
$skipRows = @(36)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) {
        continue
    }

    if ($r -eq 95) {
        # Special case: remaining count reset to 10 and start date reset.
        $ws.Cells.Item($r, 5).Value = 10
        $ws.Cells.Item($r, 6).Value = 20260301
    }
    else {
        $current = $ws.Cells.Item($r, 5).Value2
        $ws.Cells.Item($r, 5).Value = $current - 1
    }
}
